# Achievement Rename
# - Fix "Basle course of events:" -> "Basic course of events:" (cell C7 on
#   every sheet; it's a shared string reused across all four worksheets).
# - Update the active selection on each worksheet.
# - Make "Share Achievement" the active/visible sheet.

$wb = $excel.ActiveWorkbook

# Fix the typo in C7 on every worksheet (all four sheets share this string).
foreach ($ws in $wb.Worksheets) {
    $ws.Range("C7").Value = "Basic course of events:"
}

# Record Achievement: selection moves to C7.
$wsRecord = $wb.Worksheets.Item(1)
$wsRecord.Activate()
$wsRecord.Range("C7").Select()

# View Achievement: selection moves to C7.
$wsView = $wb.Worksheets.Item(2)
$wsView.Activate()
$wsView.Range("C7").Select()

# Get Point: selection moves to C7.
$wsPoint = $wb.Worksheets.Item(3)
$wsPoint.Activate()
$wsPoint.Range("C7").Select()

# Share Achievement: selection moves to C8, and this becomes the
# active/visible tab when the workbook is (re)opened.
$wsShare = $wb.Worksheets.Item(4)
$wsShare.Activate()
$wsShare.Range("C8").Select()
